# Applies the edit described by the diff:
# - Adds two new shared strings ("test1", "test2") used in A3 / A5 of "products" sheet
# - Fills A6:A47 of "products" sheet with product codes
# - Un-hides / resets formatting quirks on former rows 12-13
# - Updates selection + column widths for columns A and M

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

# Title / subtitle cells
$ws.Range("A3").Value = "test1"
$ws.Range("A5").Value = "test2"

# Product code list A6:A47
$productCodes = @(
    "AC10000-SAECO",
    "2-VC51/02",
    "2-VC51/03",
    "2-VC52/01",
    "2-VC55/01",
    "2-VC60/01",
    "2-VC64/01",
    "2-VC64/02",
    "2-VC64/03",
    "2-WM1003/01",
    "808-SCHOTT",
    "812-SCHOTT",
    "818-ORANGE-SCHOTT",
    "834-SCHOTT",
    "836-SCHOTT",
    "838-SCHOTT",
    "912-SCHOTT",
    "AB50B",
    "AC10000-SAECO",
    "AC12-KALORIK",
    "AC14-KALORIK",
    "AC2-10.5-OUT",
    "AC2-14.0-OUT",
    "AC2-7.1.-OUT",
    "AC2C10.5+",
    "AC2C14.0+",
    "AC2C3.5+",
    "AC2C5.3+",
    "AC2C7.1-IN",
    "AC2D10.5-IN",
    "AC2D14.0-IN",
    "AC2D7.1-IN",
    "AC2DH10.5-IN",
    "AC2DH14.0-IN",
    "AC2DH17.5+",
    "AC2DH7.1-IN",
    "AC2F10.5-IN",
    "AC2F14.0-IN",
    "AC2F7.0+",
    "AC2F7.1-IN",
    "AC4",
    "AC5S2.6"
)

for ($i = 0; $i -lt $productCodes.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $productCodes[$i]
}

# Rows 12/13 previously had stray formatting (custom height / hidden) from
# the old sparse layout; normalize them back to plain default rows now that
# they hold real data.
$ws.Rows.Item(13).Hidden = $false
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()

# Column A / M best-fit widths after the new content was added
$ws.Columns.Item(1).ColumnWidth = 21.335
$ws.Columns.Item(13).ColumnWidth = 5.495

# Restore the active selection
$ws.Range("E11").Select() | Out-Null

